# Atualização automática de CARLOS_BARBOSA.xlsx
#
# 1) Rename "Paineis DARQ" -> "PAINEIS DARQ"
# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3) Delete the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
[void]$wsDesarquivamentos.Delete()
$excel.DisplayAlerts = $true

Write-Host "Done: renamed sheets and removed Desarquivamentos Pendentes."
